$wb = $excel.ActiveWorkbook

# Source cell that already carries the "note" style (yellow fill / style 51)
# used for "Formatted to MM/DD/YYYY" annotations elsewhere in the workbook.
$wsStyleSrc = $wb.Worksheets.Item("ProjectTreatment")

# --- InvoicePaymentRequest sheet: add InvoicePaymentRequestDateDisplay row ---
$wsIPR = $wb.Worksheets.Item("InvoicePaymentRequest")
$wsIPR.Rows.Item(21).Insert()
$wsIPR.Range("C21").Value = "InvoicePaymentRequestDateDisplay"
$wsIPR.Range("D21").Value = "<%= invoicePaymentRequest.InvoicePaymentRequestDateDisplay %>"
$wsIPR.Range("E21").Value = "Formatted to MM/DD/YYYY"
$wsStyleSrc.Range("E10").Copy()
$wsIPR.Range("E21").PasteSpecial(-4122)

# --- Invoice sheet: add InvoiceDateDisplay row ---
$wsInv = $wb.Worksheets.Item("Invoice")
$wsInv.Rows.Item(11).Insert()
$wsInv.Range("D11").Value = "<%= invoice.InvoiceDateDisplay %>"
$wsInv.Range("C11").Value = "InvoiceDateDisplay"
$wsInv.Range("E11").Value = "Formatted to MM/DD/YYYY"
$wsStyleSrc.Range("E10").Copy()
$wsInv.Range("E11").PasteSpecial(-4122)

# The author ended their edit session with InvoicePaymentRequest as the active sheet.
$wsIPR.Activate()
